# Update the "想去人数" (number of people interested) figures on the
# sheets that contain the conference data: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1591
    $ws.Range("F3").Value = 128
    $ws.Range("F4").Value = 75
}
